$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.451.16'
$ws.Range("E2").Value = '  -0.07%  '

$ws.Range("D3").Value = '1.805.95'
$ws.Range("E3").Value = '  +0.06%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("E5").Value = '  +0.01%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '306.75'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.36%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4529'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.19%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3602'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.33%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.34'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +2.05%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07078'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.30%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8937'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +2.54%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07823'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.65%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.45'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.02%  '

$ws.Range("D14").Value = '1.876.13'
$ws.Range("E14").Value = '  +3.35%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.291'
$ws.Range("D15").ClearFormats()

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.321'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.13%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '85.69'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.88%  '

$ws.Range("E18").Value = '  -0.03%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000008495'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.71%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.006'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.08%  '

$ws.Range("D21").Value = '26.476.31'
$ws.Range("E21").Value = '  -0.08%  '

$ws.Range("E22").Value = '  +0.06%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.970'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.33%  '

$ws.Range("D24").Value = '2.091.72'
$ws.Range("E24").Value = '  +1.18%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.52'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.71%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.958'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.51%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '152.07'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.06%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.80'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.03%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.069'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +3.92%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '112.10'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.98%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.858'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.11%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08698'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.28%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.119'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.09%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.838'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +13.13%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.467'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.88%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7246'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.35%  '

$ws.Range("E37").Value = '  -0.51%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.074'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.09%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01937'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.69%  '

$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.910'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +1.56%  '

$ws.Range("B41").Value = 'Hedera'
$ws.Range("C41").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.05109'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.56%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5118'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +4.43%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.761'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.67%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1515'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -3.32%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.028'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.19%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4675'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.92%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.007'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.01%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.974'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.71%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '100.37'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.08%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.576'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.10%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05984'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.09%  '
